$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 152; this shifts existing rows 152..191 down to 153..192
$ws.Rows.Item(152).Insert()

# Populate the new row 152 with the new data record
$ws.Range("A152").Value = 9
$ws.Range("B152").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C152").Value = "Metropolitana"
$ws.Range("D152").Value = 45211
$ws.Range("D152").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E152").Value = 13
$ws.Range("F152").Value = 100112022
$ws.Range("G152").Value = "Arveja Verde"
$ws.Range("H152").Value = "Sin especificar"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 160
$ws.Range("K152").Value = 6000
$ws.Range("L152").Value = 7000
$ws.Range("M152").Value = 6500
$ws.Range("N152").Value = "$/malla 10 kilos"
$ws.Range("O152").Value = "Provincia de Melipilla"
$ws.Range("P152").Value = 650
$ws.Range("Q152").Value = 10
$ws.Range("R152").Value = "Hortaliza"
